$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 5 data - "pole wypt" with origin = 0
$ws.Cells.Item(5, 1).Value = 0          # A5
$ws.Cells.Item(5, 2).Value = 0          # B5
$ws.Cells.Item(5, 3).Value = "no"       # C5 (shared string "no")
$ws.Cells.Item(5, 6).Value = 36         # F5
$ws.Cells.Item(5, 7).Value = 36         # G5
$ws.Cells.Item(5, 8).Value = 18.321189             # H5
$ws.Cells.Item(5, 9).Value = -65.816972000000007   # I5
$ws.Cells.Item(5, 10).Value = 18.321231999999998   # J5
$ws.Cells.Item(5, 11).Value = -65.816973000000004  # K5
$ws.Cells.Item(5, 12).Value = 2028023.504987       # L5
$ws.Cells.Item(5, 13).Value = 202235.59276199999   # M5
$ws.Cells.Item(5, 14).Value = 2028028.2687830001   # N5
$ws.Cells.Item(5, 15).Value = 202235.56064800001   # O5

# Apply the same formatting used by columns L:O in existing rows (style index 1:
# numFmt "0.000000" without the centered alignment used by H:K in rows 2-4)
$ws.Range("L2").Copy()
$ws.Range("H5:O5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
